$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Spp1"
$ws.Cells.Item(2,3).Value = "Itgav"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 209.0063303333334
$ws.Cells.Item(2,8).Value = 627.018991
$ws.Cells.Item(2,9).Value = 0.6751081226665357
$ws.Cells.Item(2,10).Value = 0.6751081226665357
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 16.535604
$ws.Cells.Item(2,14).Value = 49.606812
$ws.Cells.Item(2,15).Value = 0.2120453146491552
$ws.Cells.Item(2,16).Value = 0.2120453146491552
$ws.Cells.Item(2,17).Value = 3456.045911885188
$ws.Cells.Item(2,18).Value = 31104.41320696669
$ws.Cells.Item(2,19).Value = 0.1431535142930261
$ws.Cells.Item(2,20).Value = 0.1431535142930261

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Spp1"
$ws.Cells.Item(3,3).Value = "Itgav"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 209.0063303333334
$ws.Cells.Item(3,8).Value = 627.018991
$ws.Cells.Item(3,9).Value = 0.6751081226665357
$ws.Cells.Item(3,10).Value = 0.6751081226665357
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 40.62063066666667
$ws.Cells.Item(3,14).Value = 121.861892
$ws.Cells.Item(3,15).Value = 0.5209011059384622
$ws.Cells.Item(3,16).Value = 0.5209011059384622
$ws.Cells.Item(3,17).Value = 8489.968951465666
$ws.Cells.Item(3,18).Value = 76409.72056319099
$ws.Cells.Item(3,19).Value = 0.3516645677250375
$ws.Cells.Item(3,20).Value = 0.3516645677250375

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Spp1"
$ws.Cells.Item(4,3).Value = "Itgav"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 209.0063303333334
$ws.Cells.Item(4,8).Value = 627.018991
$ws.Cells.Item(4,9).Value = 0.6751081226665357
$ws.Cells.Item(4,10).Value = 0.6751081226665357
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 20.825229
$ws.Cells.Item(4,14).Value = 62.475687
$ws.Cells.Item(4,15).Value = 0.2670535794123827
$ws.Cells.Item(4,16).Value = 0.2670535794123827
$ws.Cells.Item(4,17).Value = 4352.604691641313
$ws.Cells.Item(4,18).Value = 39173.44222477182
$ws.Cells.Item(4,19).Value = 0.1802900406484723
$ws.Cells.Item(4,20).Value = 0.1802900406484723

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Spp1"
$ws.Cells.Item(5,3).Value = "Itgav"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.401741666666666
$ws.Cells.Item(5,8).Value = 4.205225
$ws.Cells.Item(5,9).Value = 0.004527744128790482
$ws.Cells.Item(5,10).Value = 0.004527744128790482
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 16.535604
$ws.Cells.Item(5,14).Value = 49.606812
$ws.Cells.Item(5,15).Value = 0.2120453146491552
$ws.Cells.Item(5,16).Value = 0.2120453146491552
$ws.Cells.Item(5,17).Value = 23.1786451103
$ws.Cells.Item(5,18).Value = 208.6078059927
$ws.Cells.Item(5,19).Value = 0.000960086928440243
$ws.Cells.Item(5,20).Value = 0.0009600869284402431

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Spp1"
$ws.Cells.Item(6,3).Value = "Itgav"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.401741666666666
$ws.Cells.Item(6,8).Value = 4.205225
$ws.Cells.Item(6,9).Value = 0.004527744128790482
$ws.Cells.Item(6,10).Value = 0.004527744128790482
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 40.62063066666667
$ws.Cells.Item(6,14).Value = 121.861892
$ws.Cells.Item(6,15).Value = 0.5209011059384622
$ws.Cells.Item(6,16).Value = 0.5209011059384622
$ws.Cells.Item(6,17).Value = 56.93963053174444
$ws.Cells.Item(6,18).Value = 512.4566747857
$ws.Cells.Item(6,19).Value = 0.002358506924093341
$ws.Cells.Item(6,20).Value = 0.002358506924093341

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Spp1"
$ws.Cells.Item(7,3).Value = "Itgav"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.401741666666666
$ws.Cells.Item(7,8).Value = 4.205225
$ws.Cells.Item(7,9).Value = 0.004527744128790482
$ws.Cells.Item(7,10).Value = 0.004527744128790482
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 20.825229
$ws.Cells.Item(7,14).Value = 62.475687
$ws.Cells.Item(7,15).Value = 0.2670535794123827
$ws.Cells.Item(7,16).Value = 0.2670535794123827
$ws.Cells.Item(7,17).Value = 29.19159120717499
$ws.Cells.Item(7,18).Value = 262.724320864575
$ws.Cells.Item(7,19).Value = 0.001209150276256898
$ws.Cells.Item(7,20).Value = 0.001209150276256898

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Spp1"
$ws.Cells.Item(8,3).Value = "Itgav"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 99.18134533333334
$ws.Cells.Item(8,8).Value = 297.544036
$ws.Cells.Item(8,9).Value = 0.3203641332046738
$ws.Cells.Item(8,10).Value = 0.3203641332046737
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 16.535604
$ws.Cells.Item(8,14).Value = 49.606812
$ws.Cells.Item(8,15).Value = 0.2120453146491552
$ws.Cells.Item(8,16).Value = 0.2120453146491552
$ws.Cells.Item(8,17).Value = 1640.023450619248
$ws.Cells.Item(8,18).Value = 14760.21105557323
$ws.Cells.Item(8,19).Value = 0.06793171342768893
$ws.Cells.Item(8,20).Value = 0.06793171342768892

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Spp1"
$ws.Cells.Item(9,3).Value = "Itgav"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 99.18134533333334
$ws.Cells.Item(9,8).Value = 297.544036
$ws.Cells.Item(9,9).Value = 0.3203641332046738
$ws.Cells.Item(9,10).Value = 0.3203641332046737
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 40.62063066666667
$ws.Cells.Item(9,14).Value = 121.861892
$ws.Cells.Item(9,15).Value = 0.5209011059384622
$ws.Cells.Item(9,16).Value = 0.5209011059384622
$ws.Cells.Item(9,17).Value = 4028.808797808458
$ws.Cells.Item(9,18).Value = 36259.27918027611
$ws.Cells.Item(9,19).Value = 0.1668780312893314
$ws.Cells.Item(9,20).Value = 0.1668780312893313

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Spp1"
$ws.Cells.Item(10,3).Value = "Itgav"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 99.18134533333334
$ws.Cells.Item(10,8).Value = 297.544036
$ws.Cells.Item(10,9).Value = 0.3203641332046738
$ws.Cells.Item(10,10).Value = 0.3203641332046737
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 20.825229
$ws.Cells.Item(10,14).Value = 62.475687
$ws.Cells.Item(10,15).Value = 0.2670535794123827
$ws.Cells.Item(10,16).Value = 0.2670535794123827
$ws.Cells.Item(10,17).Value = 2065.474229094748
$ws.Cells.Item(10,18).Value = 18589.26806185273
$ws.Cells.Item(10,19).Value = 0.08555438848765348
$ws.Cells.Item(10,20).Value = 0.08555438848765347
